# Insert a new weekly price record at row 130 for
# "Vega Monumental Concepción - Espinaca", pushing the existing rows
# 130-155 down to 131-156.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(130).Insert()

$ws.Cells.Item(130, 1).Value = 11
$ws.Cells.Item(130, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(130, 3).Value = "Bíobío"
$ws.Cells.Item(130, 4).Value = 45275
$ws.Cells.Item(130, 5).Value = 8
$ws.Cells.Item(130, 6).Value = 100112012
$ws.Cells.Item(130, 7).Value = "Espinaca"
$ws.Cells.Item(130, 8).Value = "Sin especificar"
$ws.Cells.Item(130, 9).Value = "Primera"
$ws.Cells.Item(130, 10).Value = 150
$ws.Cells.Item(130, 11).Value = 7000
$ws.Cells.Item(130, 12).Value = 7000
$ws.Cells.Item(130, 13).Value = 7000
$ws.Cells.Item(130, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(130, 15).Value = "Región Metropolitana"
$ws.Cells.Item(130, 16).Value = 700
$ws.Cells.Item(130, 17).Value = 10
$ws.Cells.Item(130, 18).Value = "Hortaliza"
